$d = $word.ActiveDocument

# 1. Locate and remove "три нивоа " from the sentence about upstream regions,
#    turning "... домаће и три нивоа узводних регија" into
#    "... домаће и узводних регија".
$hit = $d.Content
$hit.Find.Execute("три нивоа ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$editPoint = $hit.Start
$hit.Text = ""

# 2. Word keeps a single "_GoBack" bookmark marking the most recent edit
#    location. Move it from its old spot to the point of this edit.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$goBackRange = $d.Range($editPoint, $editPoint)
$d.Bookmarks.Add("_GoBack", $goBackRange)
